# ------------------------------------------------------------------
# CompStat 94th Precinct weekly refresh: bump the report header (volume
# number + covered week) and overwrite the crime-count/percent table
# (rows 15-33) with the newly collected figures.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/number + reporting week text --------------------
$ws.Range("A8").Value = "Volume 32   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# --- Cells reverting to the "0" text placeholder (General-format ----
# shared string, matching the other N/A-style label cells) ----------
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A15").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null

# --- Row 15 ---------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 2
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -50
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 13
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = 8.333333333333
$ws.Range("L15").Value = 116.666666666667
$ws.Range("M15").Value = 1200
$ws.Range("N15").Value = 62.5

# --- Row 16 ---------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = 0
$ws.Range("E16").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 125
$ws.Range("I16").Value = 109
$ws.Range("J16").Value = 89
$ws.Range("K16").Value = 22.471910112359
$ws.Range("L16").Value = -7.627118644067
$ws.Range("M16").Value = -10.655737704918
$ws.Range("N16").Value = -80.910683012259

# --- Row 17 ---------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 137
$ws.Range("J17").Value = 128
$ws.Range("K17").Value = 7.03125
$ws.Range("L17").Value = 26.851851851851
$ws.Range("M17").Value = 98.550724637681
$ws.Range("N17").Value = -45.418326693227

# --- Row 18 ---------------------------------------------------------
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -71.428571428571
$ws.Range("J18").Value = 139
$ws.Range("K18").Value = -7.913669064748
$ws.Range("L18").Value = -21.951219512195
$ws.Range("M18").Value = -43.859649122807
$ws.Range("N18").Value = -87.96992481203

# --- Row 19 ---------------------------------------------------------
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -20.338983050847
$ws.Range("I19").Value = 567
$ws.Range("J19").Value = 645
$ws.Range("K19").Value = -12.093023255814
$ws.Range("L19").Value = -11.682242990654
$ws.Range("M19").Value = 96.875
$ws.Range("N19").Value = 64.825581395348

# --- Row 20 ---------------------------------------------------------
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 118
$ws.Range("J20").Value = 107
$ws.Range("K20").Value = 10.280373831775
$ws.Range("L20").Value = -19.727891156462
$ws.Range("M20").Value = -14.492753623188
$ws.Range("N20").Value = -86.198830409356

# --- Row 21 ---------------------------------------------------------
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -17.857142857142
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = -14.285714285714
$ws.Range("I21").Value = 1072
$ws.Range("J21").Value = 1120
$ws.Range("K21").Value = -4.285714285714
$ws.Range("L21").Value = -9.535864978902
$ws.Range("M21").Value = 26.713947990543
$ws.Range("N21").Value = -65.430506288294

# --- Row 22 ---------------------------------------------------------
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = -5.882352941176

# --- Row 23 ---------------------------------------------------------
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("G23").NumberFormat = "#,##0"
$ws.Range("H23").Value = 0
$ws.Range("H23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J23").Value = 28
$ws.Range("K23").Value = 17.857142857142

# --- Row 24 ---------------------------------------------------------
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -13.636363636363
$ws.Range("F24").Value = 56
$ws.Range("H24").Value = -28.205128205128
$ws.Range("I24").Value = 883
$ws.Range("J24").Value = 932
$ws.Range("K24").Value = -5.257510729613
$ws.Range("L24").Value = 4.869358669833
$ws.Range("M24").Value = 53.83275261324

# --- Row 25 ---------------------------------------------------------
$ws.Range("C25").Value = 7
$ws.Range("E25").Value = -46.153846153846
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = -52.727272727272
$ws.Range("I25").Value = 426
$ws.Range("J25").Value = 592
$ws.Range("K25").Value = -28.04054054054
$ws.Range("L25").Value = -8.387096774193

# --- Row 26 ---------------------------------------------------------
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 166.666666666667
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = -9.090909090909
$ws.Range("I26").Value = 262
$ws.Range("J26").Value = 259
$ws.Range("K26").Value = 1.158301158301
$ws.Range("L26").Value = 15.418502202643
$ws.Range("M26").Value = 32.994923857868

# --- Row 27 ---------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 2
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -50
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = 23.076923076923
$ws.Range("L27").Value = 100

# --- Row 28 ---------------------------------------------------------
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 46
$ws.Range("K28").Value = 17.948717948717
$ws.Range("L28").Value = 9.523809523809

# --- Row 33 ---------------------------------------------------------
$ws.Range("L33").Value = -66.666666666666
